$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# photolist.xlsx update: the photo table (rows 2-7) is rebuilt with six JV
# team entries (one row per player-less "photo set"), new Flickr / static
# Flickr hyperlinks for Photo1..Photo5, and the leftover placeholder rows
# are turned into real data rows. This also drops the now-unused
# "Drew Barklay" / "Rebels" / "Adam Hatcher" text values.
# ---------------------------------------------------------------------------

# Remove all existing hyperlinks first so re-adding them below doesn't just
# stack duplicates on top of the old D2:H3 links.
$ws.Hyperlinks.Delete()

$rows = @(
  @{ Row=2; Team="JV"; Num=1; Id=9;  Urls=@(
        "https://www.flickr.com/photos/aaronhatcher/44303421911/in/album-72157699031593671/",
        "https://www.flickr.com/photos/aaronhatcher/43597297365/in/album-72157700729583704/",
        "https://www.flickr.com/photos/aaronhatcher/29728327737/in/album-72157673392276678/",
        "https://www.flickr.com/photos/aaronhatcher/44666348631/in/album-72157673392276678/",
        "https://www.flickr.com/photos/aaronhatcher/29623981907/in/album-72157673213389378/") },
  @{ Row=3; Team="JV"; Num=2; Id=6;  Urls=@(
        "https://www.flickr.com/photos/aaronhatcher/44815854251/in/album-72157695718027700/",
        "https://www.flickr.com/photos/aaronhatcher/44921041101/in/album-72157698526626812/",
        "https://www.flickr.com/photos/aaronhatcher/29448835647/in/album-72157700819906805/",
        "https://www.flickr.com/photos/aaronhatcher/43786276754/in/album-72157701022169035/",
        "https://www.flickr.com/photos/aaronhatcher/29566593047/in/album-72157701022169035/") },
  @{ Row=4; Team="JV"; Num=3; Id=16; Urls=@(
        "https://www.flickr.com/photos/aaronhatcher/42797091140/in/album-72157673289806428/",
        "https://www.flickr.com/photos/aaronhatcher/43004907950/in/album-72157695718027700/",
        "https://www.flickr.com/photos/aaronhatcher/30795644808/in/album-72157673392276678/",
        "https://www.flickr.com/photos/aaronhatcher/44790536451/in/album-72157673602834238/",
        "https://www.flickr.com/photos/aaronhatcher/43004908230/in/album-72157695718027700/") },
  @{ Row=5; Team="JV"; Num=4; Id=18; Urls=@(
        "https://www.flickr.com/photos/aaronhatcher/44201224464/in/album-72157698526626812/",
        "https://www.flickr.com/photos/aaronhatcher/30919246858/in/album-72157673602834238/",
        "https://www.flickr.com/photos/aaronhatcher/44247376301/in/album-72157670524024857/",
        "https://www.flickr.com/photos/aaronhatcher/44255305832/in/album-72157699031593671/",
        "https://www.flickr.com/photos/aaronhatcher/43397224615/in/album-72157699031593671/") },
  @{ Row=6; Team="JV"; Num=5; Id=17; Urls=@(
        "https://farm2.staticflickr.com/1978/44346511884_7b1fff907b_o_d.jpg",
        "https://farm2.staticflickr.com/1955/31402386198_f262518aea_o_d.jpg",
        "https://farm2.staticflickr.com/1866/29624083197_49d03ec9d3_o_d.jpg",
        "https://farm2.staticflickr.com/1843/44790455141_6f20a05bbf_o_d.jpg",
        "https://farm2.staticflickr.com/1871/30634437668_9603c4f8b5_o_d.jpg") },
  @{ Row=7; Team="JV"; Num=6; Id=13; Urls=@(
        "https://www.flickr.com/photos/aaronhatcher/44666344741/in/album-72157673392276678/",
        "https://www.flickr.com/photos/aaronhatcher/44921048541/in/album-72157698526626812/",
        "https://www.flickr.com/photos/aaronhatcher/44740716672/in/album-72157673602834238/",
        "https://www.flickr.com/photos/aaronhatcher/30919341708/in/album-72157673602834238/",
        "https://www.flickr.com/photos/aaronhatcher/44740799732/in/album-72157673602834238/") }
)

$photoCols = @("D", "E", "F", "G", "H")

foreach ($r in $rows) {
  $rowNum = $r.Row

  $ws.Range("A$rowNum").Value = $r.Team
  $ws.Range("B$rowNum").Value = $r.Num
  $ws.Range("C$rowNum").Value = $r.Id

  for ($i = 0; $i -lt $photoCols.Length; $i++) {
    $addr = "$($photoCols[$i])$rowNum"
    $ws.Range($addr).Value = $r.Urls[$i]
    $ws.Hyperlinks.Add($ws.Range($addr), $r.Urls[$i]) | Out-Null
  }
}

# Rows 2-5 share one (Arial 13) look; row 6 steps to a plain-black Arial 10;
# row 7 picks up the Arial 10 "theme text" look the old row 2 used to have.
foreach ($rowNum in 2..5) {
  $rng = $ws.Range("A" + $rowNum + ":C" + $rowNum)
  $rng.Font.Name = "Arial"
  $rng.Font.Size = 13
  $rng.Font.Color = 0
  $iCell = $ws.Range("I$rowNum")
  $iCell.Font.Name = "Arial"
  $iCell.Font.Size = 13
  $iCell.Font.Color = 0
}

$row6 = $ws.Range("A6:C6")
$row6.Font.Name = "Arial"
$row6.Font.Size = 10
$row6.Font.Color = 0
$ws.Rows.Item(6).RowHeight = 17
# The old placeholder "I6" filler cell (s="2") from the before-state has no
# counterpart in the rebuilt table - drop it outright instead of leaving an
# empty styled cell behind.
$ws.Range("I6").Clear()

$row7 = $ws.Range("A7:C7")
$row7.Font.Name = "Arial"
$row7.Font.Size = 10
$row7.Font.ThemeColor = 1

# Row 3 no longer is the "tall" row (that's row 6 now) - put its height back
# to the sheet default.
$ws.Rows.Item(3).RowHeight = 16

$ws.Range("B7").Select()
